# Navigation / page lookup table addition to Feuil1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "modulo" label + value (K16 bold header, K17 the modulo constant) ---
$ws.Range("K16").Value = "modulo"
$ws.Range("K16").Font.Bold = $true

$ws.Range("K17").Value = 3

# --- "Nom Page" / "pageID" header row (bold) ---
$ws.Range("J18").Value = "Nom Page"
$ws.Range("J18").Font.Bold = $true
$ws.Range("K18").Value = "pageID"
$ws.Range("K18").Font.Bold = $true

# --- Page lookup table rows: J = page name, K = pageID, L = MOD(K,$K$17) ---
$ws.Range("J19").Value = "0 : Home Page"
$ws.Range("K19").Value = 0
$ws.Range("L19").Formula = '=MOD(K19,$K$17)'

$ws.Range("J20").Value = "1 : Ride Page"
$ws.Range("K20").Value = 1
$ws.Range("L20").Formula = '=MOD(K20,$K$17)'

$ws.Range("J21").Value = "2 : Ride Statistics Page"
$ws.Range("K21").Value = 2
$ws.Range("L21").Formula = '=MOD(K21,$K$17)'

$ws.Range("J22").Value = "3 : Compass Page"
$ws.Range("K22").Value = 3
$ws.Range("L22").Formula = '=MOD(K22,$K$17)'

$ws.Range("J23").Value = "4 : Ride Direction Page"
$ws.Range("K23").Value = 4
$ws.Range("L23").Formula = '=MOD(K23,$K$17)'

$ws.Range("J24").Value = "5 : Global Statistics Page"
$ws.Range("K24").Value = 5
$ws.Range("L24").Formula = '=MOD(K24,$K$17)'

$ws.Range("J25").Value = "6 : Go Home Page"
$ws.Range("K25").Value = 6
$ws.Range("L25").Formula = '=MOD(K25,$K$17)'

$ws.Range("J26").Value = "-1 : Init TS Page"
$ws.Range("K26").Value = -1
$ws.Range("L26").Formula = '=MOD(K26,$K$17)'

$ws.Range("J27").Value = "-2 : No Page (error)"
$ws.Range("K27").Value = -2
$ws.Range("L27").Formula = '=MOD(K27,$K$17)'

# Column J needs to be wide enough to fit the longest label (best-fit sizing)
$ws.Columns.Item(10).ColumnWidth = 21.5

# Restore the previously-selected cell to match the new working area
$ws.Range("N21").Select() | Out-Null
